# Generate Report for Handoff
# Rewrites the Overview / zh-cn / de-de sheets with the latest handoff
# status for the three source files that are currently tracked, adding
# a new row (and new hyperlinks) for the file that has just become
# "Ready for handoff".

$wb = $excel.ActiveWorkbook

$linkColor = 15570276   # OLE (BGR) encoding of RGB FF6495ED, the workbook's HyperLink font color

function Set-CellText($ws, $addr, $text) {
    $ws.Range($addr).Value = $text
}

function Style-AsHyperlink($ws, $addr) {
    $ws.Range($addr).Font.Underline = $true
    $ws.Range($addr).Font.Color = $linkColor
}

function Style-AsDateText($ws, $addr) {
    $ws.Range($addr).NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

function Add-Link($ws, $addr, $url, $display) {
    $ws.Hyperlinks.Add($ws.Range($addr), $url, [Type]::Missing, [Type]::Missing, $display) | Out-Null
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()

Set-CellText $ws1 "A1" "File Name"
Set-CellText $ws1 "B1" "zh-cn"
Set-CellText $ws1 "C1" "de-de"
Set-CellText $ws1 "D1" "Latest Handoff Date"

Set-CellText $ws1 "A2" "3a493698-11bd-44e8-8f90-052b4266a241.png"
Set-CellText $ws1 "B2" "Ready for handoff"
Set-CellText $ws1 "C2" "Ready for handoff"
Set-CellText $ws1 "D2" "2016-03-23 04:54:47"

Set-CellText $ws1 "A3" "6a9760c8-d678-4be3-aad1-a7388f2f627f.png"
Set-CellText $ws1 "B3" "Ready for handoff"
Set-CellText $ws1 "C3" "Ready for handoff"
Set-CellText $ws1 "D3" "2016-03-23 04:54:47"

Set-CellText $ws1 "A4" "95a80f8d-341a-4feb-97fe-b90923b25336.md"
Set-CellText $ws1 "B4" "Ready for handoff"
Set-CellText $ws1 "C4" "Ready for handoff"
Set-CellText $ws1 "D4" "2016-03-23 04:54:47"

Add-Link $ws1 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/926f1b677e6a3b1548907153752b30d80355cdcc/e2e/3a493698-11bd-44e8-8f90-052b4266a241.png" "3a493698-11bd-44e8-8f90-052b4266a241.png"
Add-Link $ws1 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/926f1b677e6a3b1548907153752b30d80355cdcc/e2e/6a9760c8-d678-4be3-aad1-a7388f2f627f.png" "6a9760c8-d678-4be3-aad1-a7388f2f627f.png"
Add-Link $ws1 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/926f1b677e6a3b1548907153752b30d80355cdcc/e2e/95a80f8d-341a-4feb-97fe-b90923b25336.md" "95a80f8d-341a-4feb-97fe-b90923b25336.md"

foreach ($addr in @("A2","A3","A4")) { Style-AsHyperlink $ws1 $addr }
foreach ($addr in @("D2","D3","D4")) { Style-AsDateText $ws1 $addr }

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()
$ws2.Range("F2").Clear()
$ws2.Range("G2").Clear()
$ws2.Range("F3").Clear()
$ws2.Range("G3").Clear()

Set-CellText $ws2 "A1" "Source File Name"
Set-CellText $ws2 "B1" "File Extension"
Set-CellText $ws2 "C1" "Status"
Set-CellText $ws2 "D1" "Latest Handoff File"
Set-CellText $ws2 "E1" "Latest Handoff Datetime"
Set-CellText $ws2 "F1" "Latest Target File"
Set-CellText $ws2 "G1" "Latest Handback File"
Set-CellText $ws2 "H1" "Latest Handback DateTime"
Set-CellText $ws2 "I1" "Reference Tokens"
Set-CellText $ws2 "J1" "Handoff Reason"
Set-CellText $ws2 "K1" "Dependency From"
Set-CellText $ws2 "L1" "Error Detail"

Set-CellText $ws2 "A2" "3a493698-11bd-44e8-8f90-052b4266a241.png"
Set-CellText $ws2 "B2" ".png"
Set-CellText $ws2 "C2" "Ready for handoff"
Set-CellText $ws2 "D2" "d1dc53921077e0f7d390976b7baa3c2fcb103c73.png"
Set-CellText $ws2 "E2" "2016-03-23 04:54:35"
Set-CellText $ws2 "H2" "0001-01-01 00:00:00"
Set-CellText $ws2 "J2" "IsDependency"
Set-CellText $ws2 "K2" "e2e\95a80f8d-341a-4feb-97fe-b90923b25336.md"

Set-CellText $ws2 "A3" "6a9760c8-d678-4be3-aad1-a7388f2f627f.png"
Set-CellText $ws2 "B3" ".png"
Set-CellText $ws2 "C3" "Ready for handoff"
Set-CellText $ws2 "D3" "e985d6639f04236ded45e0b898237c557dac7173.png"
Set-CellText $ws2 "E3" "2016-03-23 04:54:35"
Set-CellText $ws2 "H3" "0001-01-01 00:00:00"
Set-CellText $ws2 "J3" "IsDependency"
Set-CellText $ws2 "K3" "e2e\95a80f8d-341a-4feb-97fe-b90923b25336.md"

Set-CellText $ws2 "A4" "95a80f8d-341a-4feb-97fe-b90923b25336.md"
Set-CellText $ws2 "B4" ".md"
Set-CellText $ws2 "C4" "Ready for handoff"
Set-CellText $ws2 "D4" "95a80f8d-341a-4feb-97fe-b90923b25336.0cabf371da9956b9ff9874543b9a5419139770f9.zh-cn.xlf"
Set-CellText $ws2 "E4" "2016-03-23 04:54:35"
Set-CellText $ws2 "H4" "0001-01-01 00:00:00"
Set-CellText $ws2 "J4" "Include"

Add-Link $ws2 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/926f1b677e6a3b1548907153752b30d80355cdcc/e2e/3a493698-11bd-44e8-8f90-052b4266a241.png" "3a493698-11bd-44e8-8f90-052b4266a241.png"
Add-Link $ws2 "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd634f3c7624467043e826829a823f0d48afe980/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/d1dc53921077e0f7d390976b7baa3c2fcb103c73.png" "d1dc53921077e0f7d390976b7baa3c2fcb103c73.png"
Add-Link $ws2 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/926f1b677e6a3b1548907153752b30d80355cdcc/e2e/6a9760c8-d678-4be3-aad1-a7388f2f627f.png" "6a9760c8-d678-4be3-aad1-a7388f2f627f.png"
Add-Link $ws2 "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd634f3c7624467043e826829a823f0d48afe980/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/e985d6639f04236ded45e0b898237c557dac7173.png" "e985d6639f04236ded45e0b898237c557dac7173.png"
Add-Link $ws2 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/926f1b677e6a3b1548907153752b30d80355cdcc/e2e/95a80f8d-341a-4feb-97fe-b90923b25336.md" "95a80f8d-341a-4feb-97fe-b90923b25336.md"
Add-Link $ws2 "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd634f3c7624467043e826829a823f0d48afe980/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/95a80f8d-341a-4feb-97fe-b90923b25336.0cabf371da9956b9ff9874543b9a5419139770f9.zh-cn.xlf" "95a80f8d-341a-4feb-97fe-b90923b25336.0cabf371da9956b9ff9874543b9a5419139770f9.zh-cn.xlf"

foreach ($addr in @("A2","D2","A3","D3","A4","D4")) { Style-AsHyperlink $ws2 $addr }
foreach ($addr in @("E2","H2","E3","H3","E4","H4")) { Style-AsDateText $ws2 $addr }

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()
$ws3.Range("F2").Clear()
$ws3.Range("G2").Clear()
$ws3.Range("F3").Clear()
$ws3.Range("G3").Clear()

Set-CellText $ws3 "A1" "Source File Name"
Set-CellText $ws3 "B1" "File Extension"
Set-CellText $ws3 "C1" "Status"
Set-CellText $ws3 "D1" "Latest Handoff File"
Set-CellText $ws3 "E1" "Latest Handoff Datetime"
Set-CellText $ws3 "F1" "Latest Target File"
Set-CellText $ws3 "G1" "Latest Handback File"
Set-CellText $ws3 "H1" "Latest Handback DateTime"
Set-CellText $ws3 "I1" "Reference Tokens"
Set-CellText $ws3 "J1" "Handoff Reason"
Set-CellText $ws3 "K1" "Dependency From"
Set-CellText $ws3 "L1" "Error Detail"

Set-CellText $ws3 "A2" "3a493698-11bd-44e8-8f90-052b4266a241.png"
Set-CellText $ws3 "B2" ".png"
Set-CellText $ws3 "C2" "Ready for handoff"
Set-CellText $ws3 "D2" "d1dc53921077e0f7d390976b7baa3c2fcb103c73.png"
Set-CellText $ws3 "E2" "2016-03-23 04:54:47"
Set-CellText $ws3 "H2" "0001-01-01 00:00:00"
Set-CellText $ws3 "J2" "IsDependency"
Set-CellText $ws3 "K2" "e2e\95a80f8d-341a-4feb-97fe-b90923b25336.md"

Set-CellText $ws3 "A3" "6a9760c8-d678-4be3-aad1-a7388f2f627f.png"
Set-CellText $ws3 "B3" ".png"
Set-CellText $ws3 "C3" "Ready for handoff"
Set-CellText $ws3 "D3" "e985d6639f04236ded45e0b898237c557dac7173.png"
Set-CellText $ws3 "E3" "2016-03-23 04:54:47"
Set-CellText $ws3 "H3" "0001-01-01 00:00:00"
Set-CellText $ws3 "J3" "IsDependency"
Set-CellText $ws3 "K3" "e2e\95a80f8d-341a-4feb-97fe-b90923b25336.md"

Set-CellText $ws3 "A4" "95a80f8d-341a-4feb-97fe-b90923b25336.md"
Set-CellText $ws3 "B4" ".md"
Set-CellText $ws3 "C4" "Ready for handoff"
Set-CellText $ws3 "D4" "95a80f8d-341a-4feb-97fe-b90923b25336.0cabf371da9956b9ff9874543b9a5419139770f9.de-de.xlf"
Set-CellText $ws3 "E4" "2016-03-23 04:54:47"
Set-CellText $ws3 "H4" "0001-01-01 00:00:00"
Set-CellText $ws3 "J4" "Include"

Add-Link $ws3 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/926f1b677e6a3b1548907153752b30d80355cdcc/e2e/3a493698-11bd-44e8-8f90-052b4266a241.png" "3a493698-11bd-44e8-8f90-052b4266a241.png"
Add-Link $ws3 "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fb609572a9dee8ca05fbfb2f136bc0035554bb19/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/d1dc53921077e0f7d390976b7baa3c2fcb103c73.png" "d1dc53921077e0f7d390976b7baa3c2fcb103c73.png"
Add-Link $ws3 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/926f1b677e6a3b1548907153752b30d80355cdcc/e2e/6a9760c8-d678-4be3-aad1-a7388f2f627f.png" "6a9760c8-d678-4be3-aad1-a7388f2f627f.png"
Add-Link $ws3 "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fb609572a9dee8ca05fbfb2f136bc0035554bb19/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/e985d6639f04236ded45e0b898237c557dac7173.png" "e985d6639f04236ded45e0b898237c557dac7173.png"
Add-Link $ws3 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/926f1b677e6a3b1548907153752b30d80355cdcc/e2e/95a80f8d-341a-4feb-97fe-b90923b25336.md" "95a80f8d-341a-4feb-97fe-b90923b25336.md"
Add-Link $ws3 "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fb609572a9dee8ca05fbfb2f136bc0035554bb19/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/95a80f8d-341a-4feb-97fe-b90923b25336.0cabf371da9956b9ff9874543b9a5419139770f9.de-de.xlf" "95a80f8d-341a-4feb-97fe-b90923b25336.0cabf371da9956b9ff9874543b9a5419139770f9.de-de.xlf"

foreach ($addr in @("A2","D2","A3","D3","A4","D4")) { Style-AsHyperlink $ws3 $addr }
foreach ($addr in @("E2","H2","E3","H3","E4","H4")) { Style-AsDateText $ws3 $addr }

Write-Host "Report regenerated for handoff."
